$wb = $excel.ActiveWorkbook

# 1. CaseDetailStat sheet: add the header row (row 1) that was missing.
$ws1 = $wb.Worksheets.Item("CaseDetailStat")
$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "File Type"
$ws1.Range("C1").Value = "Association"
$ws1.Range("D1").Value = "Description"
$ws1.Range("E1").Value = "Format"
$ws1.Range("F1").Value = "Size"

# 2. CaseDetailStat_Message sheet: fix the logged Cypher query in row 28 so it
#    contains the real, resolved case id instead of the leftover 'caseid'
#    placeholder.
$ws2 = $wb.Worksheets.Item("CaseDetailStat_Message")
$ws2.Range("A28").Value = 'MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN [''NCATS-COP01CCB050022''] RETURN f.file_name AS `File Name` ,f.file_type AS `File Type`,head(labels(parent)) AS `Association`, f.file_description AS `Description`,f.file_format AS Format,((f.file_size)/1024) AS Size'
